$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings for Transportation items (A30, A32)
$ws.Range("A30").Value = "item37 [tonne*km]"
$ws.Range("A32").Value = "item38 [tonne*km]"

# Update numeric cell values
$ws.Range("F2").Value = 0.5037830471089566
$ws.Range("F3").Value = 0.5037830471089566
$ws.Range("F4").Value = 0.1356768036381773
$ws.Range("F5").Value = 0.1356768036381773
$ws.Range("F6").Value = 0.03999749251728465
$ws.Range("F7").Value = 0.03999749251728465
$ws.Range("F8").Value = 0.0001623630909792551
$ws.Range("F9").Value = 0.000268972189510174
$ws.Range("F10").Value = 0.0004313352804894291
$ws.Range("F11").Value = 0.003015040080848384
$ws.Range("F12").Value = 0.003015040080848384
$ws.Range("F13").Value = 0.003118305203617442
$ws.Range("F14").Value = 0.0001968302222476463
$ws.Range("F15").Value = 0.0004969353252158168
$ws.Range("F16").Value = 0.003812070751080905
$ws.Range("F17").Value = 0.006521531694875055
$ws.Range("F18").Value = 0.006521531694875055
$ws.Range("F19").Value = 0.2316423636194965
$ws.Range("F20").Value = 0.2316423636194965
$ws.Range("F21").Value = 0.04441464965098007
$ws.Range("F22").Value = 0.003888983380957889
$ws.Range("F23").Value = 0.04830363303193796
$ws.Range("F24").Value = 0.02747078393662985
$ws.Range("F25").Value = 0.02747078393662985
$ws.Range("E26").Value = 45438898.43131052
$ws.Range("F30").Value = 0.7769153898704085
$ws.Range("F31").Value = 0.7769153898704085
$ws.Range("C32").Value = 2678895.298605489
$ws.Range("E32").Value = 519705.6879294649
$ws.Range("F32").Value = 0.2230846101295916
$ws.Range("C33").Value = 2678895.298605489
$ws.Range("E33").Value = 519705.6879294649
$ws.Range("F33").Value = 0.2230846101295916
$ws.Range("E34").Value = 2329634.875429389
$ws.Range("B38").Value = 1824833.689404883
$ws.Range("C38").Value = 51095343.30333674
$ws.Range("D38").Value = 5.364962004765724
$ws.Range("B39").Value = 19158.22758091451
$ws.Range("C39").Value = 5076930.308942344
$ws.Range("D39").Value = 0.5330728095243165
$ws.Range("D40").Value = -0.3357271485210765
$ws.Range("B41").Value = 6565956.195095627
$ws.Range("C41").Value = -35456163.45351639
$ws.Range("D41").Value = -3.722863131256351
$ws.Range("D42").Value = -0.1301955148752136
$ws.Range("D43").Value = -0.1203654972634857
$ws.Range("B44").Value = 565659.3778841281
$ws.Range("C44").Value = -3054560.640574292
$ws.Range("D44").Value = -0.3207259354467191
$ws.Range("D45").Value = -0.2681575869271945
$ws.Range("C46").Value = 9523896.582668895
$ws.Range("C50").Value = 51864.96
$ws.Range("C51").Value = 51864.96
